$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen the first two columns (matches the saved column metadata)
$ws.Columns.Item(1).ColumnWidth = 39.5
$ws.Columns.Item(2).ColumnWidth = 41.333333333333336

# Row 2 - duplicate of the applicant's test record
$ws.Range("A2").Value = "Test Automation"
$ws.Range("B2").Value = "28/10/2019"
$ws.Range("C2").Value = "Y000-102-01"

# Row 3
$ws.Range("A3").Value = "Test Automation"
$ws.Range("B3").Value = "28/10/2019"
$ws.Range("C3").Value = "Y000-102-01"

# Row 4
$ws.Range("A4").Value = "Test Automation"
$ws.Range("B4").Value = "28/10/2019"
$ws.Range("C4").Value = "Y000-102-01"

# Row 5
$ws.Range("A5").Value = "Test Automation"
$ws.Range("B5").Value = "24/11/2019"
$ws.Range("C5").Value = "M000-097-01"

# Row 6 - DOB stored as a real date value with builtin date format (numFmtId 14)
$ws.Range("B6").NumberFormat = "mm-dd-yy"
$ws.Range("A6").Value = "Test Automation"
$ws.Range("B6").Value = (Get-Date -Year 2017 -Month 7 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C6").Value = "T000-105-01"

# Leave the selection on the last-touched cell (closest reachable approximation
# of the saved multi-area selection "A2 B4 B3 C4" with active cell C4)
$ws.Range("C4").Select()
